$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.224.91'
$ws.Range("E2").Value = '  +1.19%  '

$ws.Range("D3").Value = '1.835.18'
$ws.Range("E3").Value = '  +1.04%  '

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.014'
$c.ClearFormats()
$ws.Range("E4").Value = '  +1.31%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '313.93'
$c.ClearFormats()
$ws.Range("E5").Value = '  +1.31%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '1.012'
$c.ClearFormats()
$ws.Range("E6").Value = '  +1.11%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.4715'
$c.ClearFormats()
$ws.Range("E7").Value = '  +1.08%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3686'
$c.ClearFormats()
$ws.Range("E8").Value = '  -0.34%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.07422'
$c.ClearFormats()
$ws.Range("E9").Value = '  +0.66%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.8827'
$c.ClearFormats()
$ws.Range("E10").Value = '  +1.28%  '

$ws.Range("E11").Value = '  +0.02%  '

$ws.Range("D12").Value = '1.840.63'
$ws.Range("E12").Value = '  +1.38%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.07337'
$c.ClearFormats()
$ws.Range("E13").Value = '  +3.62%  '

$ws.Range("E14").Value = '  +1.91%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '93.11'
$c.ClearFormats()
$ws.Range("E15").Value = '  +0.74%  '

$ws.Range("E16").Value = '  +1.02%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.014'
$c.ClearFormats()
$ws.Range("E17").Value = '  +1.24%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.000008803'
$c.ClearFormats()
$ws.Range("E18").Value = '  +0.97%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.80'
$c.ClearFormats()
$ws.Range("E20").Value = '  +0.28%  '

$ws.Range("D21").Value = '27.244.66'
$ws.Range("E21").Value = '  +1.15%  '

$ws.Range("E22").Value = '  -0.73%  '

$ws.Range("E23").Value = '  +1.24%  '

$ws.Range("D24").Value = '2.069.06'
$ws.Range("E24").Value = '  +1.40%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '1.904'
$c.ClearFormats()
$ws.Range("E25").Value = '  +0.16%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '153.15'
$c.ClearFormats()
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("E27").Value = '  +1.16%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.168'
$c.ClearFormats()
$ws.Range("E28").Value = '  -1.67%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.278'
$c.ClearFormats()
$ws.Range("E29").Value = '  -0.73%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '117.82'
$c.ClearFormats()
$ws.Range("E30").Value = '  +1.89%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.08929'
$c.ClearFormats()
$ws.Range("E31").Value = '  +0.05%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.7585'
$c.ClearFormats()
$ws.Range("E32").Value = '  -1.30%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.172'
$c.ClearFormats()
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("E34").Value = '  +1.35%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '2.947'
$c.ClearFormats()
$ws.Range("E35").Value = '  +1.04%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.012'
$c.ClearFormats()
$ws.Range("E36").Value = '  +1.12%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.104'
$c.ClearFormats()
$ws.Range("E37").Value = '  +1.01%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.05335'
$c.ClearFormats()
$ws.Range("E38").Value = '  +1.06%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.01958'
$c.ClearFormats()
$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("E40").Value = '  +1.17%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '7.335'
$c.ClearFormats()
$ws.Range("E41").Value = '  +1.06%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.403'
$c.ClearFormats()
$ws.Range("E42").Value = '  +1.16%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.5342'
$c.ClearFormats()
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("E44").Value = '  +0.01%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.539'
$c.ClearFormats()
$ws.Range("E45").Value = '  +1.00%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.4956'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.24%  '

$ws.Range("E47").Value = '  +1.01%  '

$ws.Range("E48").Value = '  +1.19%  '

$ws.Range("E49").Value = '  +0.01%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '103.99'
$c.ClearFormats()
$ws.Range("E50").Value = '  +1.16%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.06324'
$c.ClearFormats()
$ws.Range("E51").Value = '  +0.50%  '

